$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F6").Value = 299
$ws.Range("F7").Value = 13316
$ws.Range("F9").Value = 342
$ws.Range("G9").Value = 40
$ws.Range("F10").Value = 5322
$ws.Range("F16").Value = 185
$ws.Range("F19").Value = 87
$ws.Range("F21").Value = 3780
$ws.Range("F22").Value = 116
$ws.Range("F24").Value = 5021
$ws.Range("F26").Value = 2005
$ws.Range("F28").Value = 302
$ws.Range("F29").Value = 7359
$ws.Range("F33").Value = 2102
$ws.Range("F34").Value = 1317
$ws.Range("F35").Value = 136
$ws.Range("F36").Value = 1138
$ws.Range("F37").Value = 11
$ws.Range("F38").Value = 242
$ws.Range("F40").Value = 6
$ws.Range("F41").Value = 1164
$ws.Range("F43").Value = 17
$ws.Range("F44").Value = 159
$ws.Range("F45").Value = 1276
$ws.Range("F46").Value = 1935
$ws.Range("F47").Value = 95
$ws.Range("F49").Value = 1195

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 11

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 515
$ws.Range("F3").Value = 689

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F6").Value = 515
$ws.Range("F7").Value = 689
$ws.Range("F8").Value = 299
$ws.Range("F9").Value = 13316
$ws.Range("F10").Value = 342
$ws.Range("G10").Value = 40
$ws.Range("F11").Value = 5322
$ws.Range("F14").Value = 185
$ws.Range("F16").Value = 87
$ws.Range("F17").Value = 11
$ws.Range("F20").Value = 3780
$ws.Range("F22").Value = 116
$ws.Range("F23").Value = 5021
$ws.Range("F25").Value = 2005
$ws.Range("F27").Value = 302
$ws.Range("F28").Value = 7359
$ws.Range("F32").Value = 2102
$ws.Range("F33").Value = 1317
$ws.Range("F34").Value = 136
$ws.Range("F35").Value = 1138
$ws.Range("F36").Value = 11
$ws.Range("F37").Value = 242
$ws.Range("F39").Value = 6
$ws.Range("F40").Value = 1164
$ws.Range("F42").Value = 17
$ws.Range("F43").Value = 159
$ws.Range("F45").Value = 1276
$ws.Range("F46").Value = 1935
$ws.Range("F47").Value = 95
$ws.Range("F49").Value = 1195
